$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-76 down to 41-77.
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with its data.
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C40").Value = "Arica y Parinacota"
$ws.Range("D40").Value = 44629
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 100112038
$ws.Range("G40").Value = "Cebollín baby"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 250
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = 1750
$ws.Range("N40").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O40").Value = "Región de Arica y Parinacota"
$ws.Range("P40").Value = 875
$ws.Range("Q40").Value = 2
$ws.Range("R40").Value = "Hortaliza"
